$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1073.6666
$ws.Range("I19").Value = 990.6667
$ws.Range("J19").Value = 1156.6666
$ws.Range("K19").Value = 990.6667
$ws.Range("L19").Value = 1156.6666
$ws.Range("M19").Value = -815.6667
$ws.Range("N19").Value = -1506.6666
$ws.Range("H103").Value = 477.2
$ws.Range("I103").Value = 532.4
$ws.Range("K103").Value = 1597.2
$ws.Range("M103").Value = -1011.2
$ws.Range("H112").Value = 2240.4666
$ws.Range("J112").Value = 2442.3333
$ws.Range("L112").Value = 7326.999899999999
$ws.Range("N112").Value = -9542.999899999999
$ws.Range("H132").Value = 4040.9153
$ws.Range("I132").Value = 4142.964
$ws.Range("K132").Value = 12428.892
$ws.Range("M132").Value = -9898.892
$ws.Range("H138").Value = 5649.57
$ws.Range("J138").Value = 5844.5957
$ws.Range("L138").Value = 17533.7871
$ws.Range("N138").Value = -27813.7871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11771.429
$ws.Range("I32").Value = 11316.129
$ws.Range("K32").Value = 11316.129
$ws.Range("M32").Value = -11029.129
$ws.Range("H88").Value = 12123.5
$ws.Range("I88").Value = 1699
$ws.Range("J88").Value = 13071.182
$ws.Range("K88").Value = 1699
$ws.Range("L88").Value = 13071.182
$ws.Range("M88").Value = -1293
$ws.Range("N88").Value = -13883.182
$ws.Range("H91").Value = 12123.5
$ws.Range("I91").Value = 1699
$ws.Range("J91").Value = 13071.182
$ws.Range("K91").Value = 1699
$ws.Range("L91").Value = 13071.182
$ws.Range("M91").Value = -295
$ws.Range("N91").Value = -15879.182
$ws.Range("H97").Value = 3476.4482
$ws.Range("I97").Value = 1734.5714
$ws.Range("J97").Value = 5102.2
$ws.Range("K97").Value = 1734.5714
$ws.Range("L97").Value = 5102.2
$ws.Range("M97").Value = -1238.5714
$ws.Range("N97").Value = -6094.2
$ws.Range("H110").Value = 1069.2632
$ws.Range("J110").Value = 2625
$ws.Range("L110").Value = 2625
$ws.Range("N110").Value = -6715
$ws.Range("H132").Value = 3435.4814
$ws.Range("I132").Value = 2647.2195
$ws.Range("J132").Value = 5921.5386
$ws.Range("K132").Value = 7941.6585
$ws.Range("L132").Value = 17764.6158
$ws.Range("M132").Value = -5411.6585
$ws.Range("N132").Value = -22824.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3610.8333
$ws.Range("I20").Value = 4341.5
$ws.Range("J20").Value = 2149.5
$ws.Range("K20").Value = 4341.5
$ws.Range("L20").Value = 2149.5
$ws.Range("M20").Value = -4094.5
$ws.Range("N20").Value = -2643.5
$ws.Range("H105").Value = 2673.2856
$ws.Range("I105").Value = 1766.4546
$ws.Range("J105").Value = 5998.3335
$ws.Range("K105").Value = 1766.4546
$ws.Range("L105").Value = 5998.3335
$ws.Range("M105").Value = -19.45460000000003
$ws.Range("N105").Value = -9492.333500000001
$ws.Range("H134").Value = 6399.586
$ws.Range("I134").Value = 6448.4443
$ws.Range("J134").Value = 6234.6875
$ws.Range("K134").Value = 19345.3329
$ws.Range("L134").Value = 18704.0625
$ws.Range("M134").Value = -16810.3329
$ws.Range("N134").Value = -23774.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3832.125
$ws.Range("I58").Value = 3439.4614
$ws.Range("K58").Value = 3439.4614
$ws.Range("M58").Value = -3236.4614
$ws.Range("H62").Value = 3457.5
$ws.Range("I62").Value = 2821
$ws.Range("J62").Value = 4518.3335
$ws.Range("K62").Value = 2821
$ws.Range("L62").Value = 4518.3335
$ws.Range("M62").Value = -2197
$ws.Range("N62").Value = -5766.3335
$ws.Range("H65").Value = 3457.5
$ws.Range("I65").Value = 2821
$ws.Range("J65").Value = 4518.3335
$ws.Range("K65").Value = 14105
$ws.Range("L65").Value = 22591.6675
$ws.Range("M65").Value = -10985
$ws.Range("N65").Value = -28831.6675
$ws.Range("H136").Value = 3832.125
$ws.Range("I136").Value = 3439.4614
$ws.Range("K136").Value = 10318.3842
$ws.Range("M136").Value = -7768.3842
$ws.Range("H141").Value = 885649.7
$ws.Range("J141").Value = 1303474.5
$ws.Range("L141").Value = 1303474.5
$ws.Range("N141").Value = -1313834.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 974
$ws.Range("J44").Value = 750
$ws.Range("L44").Value = 2250
$ws.Range("N44").Value = -3046

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H70").Value = 8727
$ws.Range("I70").Value = 9794.23
$ws.Range("K70").Value = 9794.23
$ws.Range("M70").Value = -9524.23
$ws.Range("H73").Value = 8727
$ws.Range("I73").Value = 9794.23
$ws.Range("K73").Value = 9794.23
$ws.Range("M73").Value = -8858.23
$ws.Range("H102").Value = 13983
$ws.Range("I102").Value = 15722.474
$ws.Range("K102").Value = 15722.474
$ws.Range("M102").Value = -14100.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1381.091
$ws.Range("I22").Value = 965
$ws.Range("K22").Value = 965
$ws.Range("M22").Value = -670
$ws.Range("H27").Value = 1381.091
$ws.Range("I27").Value = 965
$ws.Range("K27").Value = 965
$ws.Range("M27").Value = -858
$ws.Range("H68").Value = 2316.3125
$ws.Range("I68").Value = 2255.1667
$ws.Range("J68").Value = 2499.75
$ws.Range("K68").Value = 2255.1667
$ws.Range("L68").Value = 2499.75
$ws.Range("M68").Value = -1506.1667
$ws.Range("N68").Value = -3997.75
$ws.Range("H71").Value = 2316.3125
$ws.Range("I71").Value = 2255.1667
$ws.Range("J71").Value = 2499.75
$ws.Range("K71").Value = 11275.8335
$ws.Range("L71").Value = 12498.75
$ws.Range("M71").Value = -7531.833500000001
$ws.Range("N71").Value = -19986.75
$ws.Range("H93").Value = 1465
$ws.Range("J93").Value = 1389.4
$ws.Range("L93").Value = 1389.4
$ws.Range("N93").Value = -3885.4
$ws.Range("H100").Value = 2103.9092
$ws.Range("I100").Value = 1714.4
$ws.Range("K100").Value = 1714.4
$ws.Range("M100").Value = -1173.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 27500
$ws.Range("I15").Value = 5000
$ws.Range("J15").Value = 50000
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 50000
$ws.Range("M15").Value = -4712
$ws.Range("N15").Value = -50576
$ws.Range("H54").Value = 51610.332
$ws.Range("I54").Value = 49932.4
$ws.Range("J54").Value = 60000
$ws.Range("K54").Value = 49932.4
$ws.Range("L54").Value = 60000
$ws.Range("M54").Value = -49412.4
$ws.Range("N54").Value = -61040
$ws.Range("H81").Value = 7859.3335
$ws.Range("I81").Value = 11984.286
$ws.Range("J81").Value = 4250
$ws.Range("K81").Value = 23968.572
$ws.Range("L81").Value = 8500
$ws.Range("M81").Value = -22907.572
$ws.Range("N81").Value = -10622
$ws.Range("H84").Value = 7859.3335
$ws.Range("I84").Value = 11984.286
$ws.Range("J84").Value = 4250
$ws.Range("K84").Value = 119842.86
$ws.Range("L84").Value = 42500
$ws.Range("M84").Value = -114538.86
$ws.Range("N84").Value = -53108
